$d = $word.ActiveDocument

# --- Problem 1 ---
$d.Content.Find.Execute(
    "Lily has three balloons. Sam gives her four more. How many balloons does Lily have now?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Barnaby the badger found eight shiny pebbles. He gave three pebbles to his friend Penelope. How many pebbles does Barnaby have now?",
    2)

$d.Content.Find.Execute(
    "Answer: Lily has ____________ balloons now.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answer: Barnaby has ________________________________________ pebbles now.",
    2)

# --- Problem 2 ---
$d.Content.Find.Execute(
    "David had eleven stickers. He gave five stickers to his friend. How many stickers does David have left?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Flora the fairy baked seven cupcakes. Her friend Fizz ate five of them. How many cupcakes are left?",
    2)

$d.Content.Find.Execute(
    "Answer: David has ____________ stickers left.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answer: There are ________________________________________ cupcakes left.",
    2)

# --- Problem 3 ---
$d.Content.Find.Execute(
    "Sarah saw eight ladybugs in the garden. Two more ladybugs landed nearby. How many ladybugs are there in total?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Gregory the gnome had four mushrooms. He found six more mushrooms in the forest. How many mushrooms does Gregory have in total?",
    2)

$d.Content.Find.Execute(
    "Answer: There are ____________ ladybugs in total.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answer: Gregory has ________________________________________ mushrooms in total.",
    2)

# --- Problem 4 ---
$d.Content.Find.Execute(
    "Ben had nine crayons. He found six more under his bed. How many crayons does Ben have now?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Princess Petunia had nine sparkly wands. She bought four more wands from the magic shop. How many wands does she have altogether?",
    2)

$d.Content.Find.Execute(
    "Answer: Ben now has ____________ crayons.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answer: Princess Petunia has ________________________________________ wands altogether.",
    2)

# --- Problem 5 ---
$d.Content.Find.Execute(
    "Chloe baked fifteen cookies. Her brother ate three of them. How many cookies are left?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Professor Bumble had twelve test tubes. Two of the test tubes broke during an experiment. How many test tubes does Professor Bumble have left?",
    2)

$d.Content.Find.Execute(
    "Answer: There are ____________ cookies left.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answer: Professor Bumble has ________________________________________ test tubes left.",
    2)

# --- Problem 6 ---
$d.Content.Find.Execute(
    "Peter has seven toy cars. Lisa has eight toy cars. How many toy cars do they have altogether?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Silly Sally had five bouncy balls. Then, her dog Buttons gave her three more. How many bouncy balls does Sally have now?",
    2)

$d.Content.Find.Execute(
    "Answer: They have ____________ toy cars altogether.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answer: Silly Sally now has ________________________________________ bouncy balls.",
    2)

# --- Problem 7 ---
$d.Content.Find.Execute(
    "Emily had twelve candies. She ate four candies. How many candies does Emily have left?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Captain Calico had thirteen gold doubloons. He spent one doubloon on a parrot. How many doubloons does he have left?",
    2)

$d.Content.Find.Execute(
    "Answer: Emily has ____________ candies left.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answer: Captain Calico has ________________________________________ doubloons left.",
    2)

# --- Shared "Number sentence:" blanks (identical across all 7 paragraphs) ---
$d.Content.Find.Execute(
    "Number sentence: __________________________________________________",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Number sentence: ______________________________________________________________________",
    2)
